# Applies the latest cryptos-list price/volume refresh to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the text must be forced
# to stay text (i.e. it would otherwise be auto-parsed as a number by Excel,
# which would strip things like trailing zeros or the thousands-dot formatting).
$updates = @(
    @{ Cell = "D2"; Value = "29.185.89"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +0.16%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "1.855.91"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +0.14%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.16%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "0.7052"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.00%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "237.38"; ForceText = $true }
    @{ Cell = "E6"; Value = "  -0.62%  "; ForceText = $false }
    @{ Cell = "E7"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.07940"; ForceText = $true }
    @{ Cell = "E8"; Value = "  +4.01%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.3015"; ForceText = $true }
    @{ Cell = "E9"; Value = "  -1.30%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "23.35"; ForceText = $true }
    @{ Cell = "E10"; Value = "  -0.55%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.08168"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +0.37%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "1.869.45"; ForceText = $false }
    @{ Cell = "E12"; Value = "  +3.87%  "; ForceText = $false }
    @{ Cell = "E13"; Value = "  -1.84%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "0.7002"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -3.80%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "89.47"; ForceText = $true }
    @{ Cell = "E15"; Value = "  -0.10%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "29.169.85"; ForceText = $false }
    @{ Cell = "E16"; Value = "  +0.22%  "; ForceText = $false }
    @{ Cell = "B17"; Value = "Uniswap"; ForceText = $false }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; ForceText = $false }
    @{ Cell = "D17"; Value = "5.780"; ForceText = $true }
    @{ Cell = "E17"; Value = "  -0.35%  "; ForceText = $false }
    @{ Cell = "B18"; Value = "ShibaInu"; ForceText = $false }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; ForceText = $false }
    @{ Cell = "D18"; Value = "0.000007843"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +0.81%  "; ForceText = $false }
    @{ Cell = "E19"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "235.45"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -1.03%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "1.000"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "2.094.47"; ForceText = $false }
    @{ Cell = "E22"; Value = "  -0.05%  "; ForceText = $false }
    @{ Cell = "E23"; Value = "  +0.18%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "7.341"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -3.43%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "161.86"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +0.50%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "8.885"; ForceText = $true }
    @{ Cell = "E26"; Value = "  -1.68%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "0.1426"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -2.15%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "17.98"; ForceText = $true }
    @{ Cell = "E28"; Value = "  -0.89%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "1.923"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -3.44%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "1.428"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +1.90%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "1.478"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -1.07%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "4.354"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -3.37%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "4.010"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -0.05%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "0.05188"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -0.65%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "1.157"; ForceText = $true }
    @{ Cell = "E35"; Value = "  -3.00%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "0.6981"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -1.14%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "1.001"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -2.63%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "2.669"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -0.04%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "0.01835"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -1.78%  "; ForceText = $false }
    @{ Cell = "E40"; Value = "  +1.32%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "0.9298"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -1.08%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "1.119.58"; ForceText = $false }
    @{ Cell = "E42"; Value = "  +4.04%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "0.4232"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -1.41%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "5.830"; ForceText = $true }
    @{ Cell = "E44"; Value = "  -3.31%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "69.37"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -1.42%  "; ForceText = $false }
    @{ Cell = "E46"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "102.49"; ForceText = $true }
    @{ Cell = "E47"; Value = "  -0.55%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "1.750"; ForceText = $true }
    @{ Cell = "E48"; Value = "  -1.74%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "1.991.66"; ForceText = $false }
    @{ Cell = "E49"; Value = "  +0.93%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "9.077"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -1.52%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "0.05929"; ForceText = $true }
    @{ Cell = "E51"; Value = "  +1.02%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Prefix with an apostrophe (as typing in the Excel UI would) so the
        # numeric-looking text is stored as a literal string, then restore the
        # default "Normal" cell style so no stray number format sticks around.
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
